# Auto-generated edit script applying the cryptos.xlsx diff
# Updates Price (D) and Volume(1h) (E) columns, plus one row swap (B26:E27)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Cell, $Text)
    $rng = $ws.Range($Cell)
    # Force text type so numeric-looking strings (e.g. "303.38") are not
    # auto-coerced into numbers by Excel - matches the source workbook,
    # which stores every one of these cells as inline text.
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    # Drop back to the default (unstyled) cell format - none of these
    # cells carry an explicit style in the source file.
    $rng.ClearFormats()
}

Set-TextCell 'D2' '41.742.11'
Set-TextCell 'E2' '  +1.40%  '
Set-TextCell 'D3' '2.265.46'
Set-TextCell 'E3' '  +0.86%  '
Set-TextCell 'E4' '  -0.07%  '
Set-TextCell 'D5' '303.38'
Set-TextCell 'D6' '92.03'
Set-TextCell 'E6' '  +1.23%  '
Set-TextCell 'E7' '  +2.03%  '
Set-TextCell 'E8' '  -0.11%  '
Set-TextCell 'E9' '  -0.04%  '
Set-TextCell 'D10' '32.39'
Set-TextCell 'E11' '  -1.09%  '
Set-TextCell 'E12' '  +0.45%  '
Set-TextCell 'E13' '  -1.35%  '
Set-TextCell 'E14' '  +1.16%  '
Set-TextCell 'D15' '2.616.19'
Set-TextCell 'E15' '  +0.77%  '
Set-TextCell 'D16' '14.23'
Set-TextCell 'E16' '  +1.11%  '
Set-TextCell 'D17' '2.275.67'
Set-TextCell 'E17' '  +2.93%  '
Set-TextCell 'D18' '0.768'
Set-TextCell 'D19' '41.640.85'
Set-TextCell 'E19' '  +1.31%  '
Set-TextCell 'D20' '12.48'
Set-TextCell 'E20' '  +5.50%  '
Set-TextCell 'E21' '  +0.45%  '
Set-TextCell 'D22' '5.94'
Set-TextCell 'E22' '  +1.67%  '
Set-TextCell 'D23' '67.10'
Set-TextCell 'E23' '  +0.42%  '
Set-TextCell 'D24' '239.63'
Set-TextCell 'E24' '  -0.25%  '
Set-TextCell 'D25' '2.58'
Set-TextCell 'E25' '  +1.08%  '
Set-TextCell 'B26' 'Dai'
Set-TextCell 'C26' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell 'D26' '1.00'
Set-TextCell 'E26' '  -0.15%  '
Set-TextCell 'B27' 'ImmutableX'
Set-TextCell 'C27' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 'D27' '1.92'
Set-TextCell 'E27' '  +4.28%  '
Set-TextCell 'D28' '23.91'
Set-TextCell 'E28' '  +0.85%  '
Set-TextCell 'D29' '9.54'
Set-TextCell 'E29' '  -0.17%  '
Set-TextCell 'D30' '2.08'
Set-TextCell 'E30' '  -0.84%  '
Set-TextCell 'D31' '35.38'
Set-TextCell 'E31' '  +7.09%  '
Set-TextCell 'D32' '159.90'
Set-TextCell 'E32' '  +0.57%  '
Set-TextCell 'D33' '5.25'
Set-TextCell 'E33' '  +1.71%  '
Set-TextCell 'D34' '0.999'
Set-TextCell 'E34' '  -0.14%  '
Set-TextCell 'D35' '0.0742'
Set-TextCell 'E35' '  +1.78%  '
Set-TextCell 'E36' '  -0.08%  '
Set-TextCell 'D37' '16.91'
Set-TextCell 'E37' '  +2.72%  '
Set-TextCell 'E38' '  +0.48%  '
Set-TextCell 'E39' '  +1.68%  '
Set-TextCell 'E40' '  +1.18%  '
Set-TextCell 'D41' '1.78'
Set-TextCell 'E41' '  +0.86%  '
Set-TextCell 'E42' '  -0.09%  '
Set-TextCell 'D43' '2.014.11'
Set-TextCell 'E43' '  -2.88%  '
Set-TextCell 'D44' '19.28'
Set-TextCell 'E44' '  -4.03%  '
Set-TextCell 'E45' '  +1.19%  '
Set-TextCell 'D46' '10.35'
Set-TextCell 'E46' '  +1.44%  '
Set-TextCell 'D47' '2.11'
Set-TextCell 'E47' '  +5.05%  '
Set-TextCell 'E48' '  -1.57%  '
Set-TextCell 'E49' '  +1.95%  '
Set-TextCell 'E50' '  +0.83%  '
Set-TextCell 'D51' '52.19'
Set-TextCell 'E51' '  +3.02%  '
